$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.102.33"
$ws.Range("E2").Value = "  +3.38%  "
$ws.Range("D3").Value = "2.447.66"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'322.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").Value = "'104.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.77%  "
$ws.Range("D10").Value = "'36.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").Value = "'0.0806"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "'18.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "2.834.72"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "2.461.21"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "45.929.77"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "'12.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "'6.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  +2.88%  "
$ws.Range("D22").Value = "'71.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("D23").Value = "'247.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").Value = "'2.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").Value = "'26.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.99%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'2.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.34%  "
$ws.Range("D29").Value = "'9.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "'33.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").Value = "'49.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  +6.17%  "
$ws.Range("D33").Value = "'20.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'0.0762"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "'127.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("E41").Value = "  +5.90%  "
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "'21.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").Value = "1.959.52"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("E48").Value = "  +9.55%  "
$ws.Range("D49").Value = "'9.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.59%  "
$ws.Range("D50").Value = "'77.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.23%  "
$ws.Range("D51").Value = "'4.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.54%  "
